$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "Runtime in s"
$ws.Range("H2").Value = 16
$ws.Range("H3").Value = 186
$ws.Range("H5").Value = 29
$ws.Range("H6").Value = 89
$ws.Range("H1").Style = "Comma"
$ws.Range("H2").Style = "Comma"
$ws.Range("H3").Style = "Comma"
$ws.Range("H5").Style = "Comma"
$ws.Range("H6").Style = "Comma"
$ws.Range("H1").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws.Range("H2").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws.Range("H3").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws.Range("H5").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws.Range("H6").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
Write-Host "done"
